# Apply the "simulator full-month coverage, persist logs, fix employees" edit.
#
# - Fix employee/client names that the simulator had wrong.
# - Persist the simulator's computed Rate/Total figures (previously 0/0
#   placeholders) onto both the "Weekly Timesheet" sheet and the
#   "Jason Schema" mirror sheet.

$wb = $excel.ActiveWorkbook

$wsTime   = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Fix employee data -----------------------------------------------------

# Client names on the "Weekly Timesheet" sheet.
$wsTime.Range("B2").Value = "McClure"
$wsTime.Range("B3").Value = "Evans"

# Mirrored client names + employee id on the "Jason Schema" sheet.
$wsSchema.Range("D2").Value = "McClure"
$wsSchema.Range("D3").Value = "Evans"
$wsSchema.Range("B2").Value = "emp_4nlnrvy7"
$wsSchema.Range("B3").Value = "emp_4nlnrvy7"

# --- Persist simulator Rate / Total figures --------------------------------

# "Weekly Timesheet" sheet: daily rows.
$wsTime.Range("E2").Value = 95
$wsTime.Range("F2").Value = 1900
$wsTime.Range("E3").Value = 95
$wsTime.Range("F3").Value = 1900

# "Weekly Timesheet" sheet: subtotal rows.
$wsTime.Range("F5").Value = 3800
$wsTime.Range("F8").Value = 3800
$wsTime.Range("F10").Value = 3800

# "Jason Schema" sheet: mirrored Rate / Total columns.
$wsSchema.Range("F2").Value = 95
$wsSchema.Range("G2").Value = 1900
$wsSchema.Range("F3").Value = 95
$wsSchema.Range("G3").Value = 1900
